# Generate Report for Handoff
# Adds a new localization-status row (file ac9b1417-61ed-4e9d-9460-3b986c173417.md)
# to each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e6d0bfe9517d0f1dd06ebce4f84d0a9094908b58/e2e/"
$newFile = "ac9b1417-61ed-4e9d-9460-3b986c173417.md"
$newFileUrl = $repoBase + $newFile

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-31 00:42:35"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", "e2e\" + $newFile)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "ac9b1417-61ed-4e9d-9460-3b986c173417.d5e72f5512c9f7d61d8086a6de28ce4d78648fa1.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-31 00:42:30"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newFileUrl, "", "", $newFile)

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "ac9b1417-61ed-4e9d-9460-3b986c173417.d5e72f5512c9f7d61d8086a6de28ce4d78648fa1.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-31 00:42:35"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newFileUrl, "", "", $newFile)

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

Write-Output "Added handoff row for $newFile to Overview, zh-cn, de-de sheets"
